{"js": "const body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\n// \"{{ Name_2 \" + _GoBack bookmark + \"}}\" (two runs split by a stray\n// bookmark) collapses into a single run \"{{ Name_2 }}\" and the bookmark\n// that used to sit between the runs is removed.\nconst name2Results = body.search(\"{{ Name_2 }}\", { matchCase: true });\nawait context.sync();\n\nif (name2Results.items.length > 0) {\n  const name2Ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' +\n              '<w:p>' +\n                '<w:r w:rsidRPr=\"00EA0001\">' +\n                  '<w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n                  '<w:t>{{ Name_2 }}</w:t>' +\n                '</w:r>' +\n              '</w:p>' +\n            '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>';\n  name2Results.items[0].insertOoxml(name2Ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2 -------------------------------------------------------------\n// \"Firma del Dise\u00f1ador\" becomes \"Firma de\" + a new _GoBack bookmark +\n// \" Dise\u00f1o\" (still bold) - i.e. the bookmark now sits in the middle of the\n// (renamed) signature label instead of inside the Name_2 placeholder.\nconst designerResults = body.search(\"Firma del Dise\u00f1ador\", { matchCase: true });\nawait context.sync();\n\nif (designerResults.items.length > 0) {\n  const designerOoxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' +\n              '<w:p>' +\n                '<w:r>' +\n                  '<w:rPr><w:b/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n                  '<w:t>Firma de</w:t>' +\n                '</w:r>' +\n                '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n                '<w:bookmarkEnd w:id=\"0\"/>' +\n                '<w:r>' +\n                  '<w:rPr><w:b/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n                  '<w:t xml:space=\"preserve\"> Dise\\u00f1o</w:t>' +\n                '</w:r>' +\n              '</w:p>' +\n            '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>';\n  designerResults.items[0].insertOoxml(designerOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: \"{{ Name_2 \" + _GoBack bookmark + \"}}\" (two runs) -> single run \"{{ Name_2 }}\" ---\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Text = \"{{ Name_2 }}\"\n$rng1.Find.Execute() | Out-Null\nif ($rng1.Find.Found) {\n    # Force a genuine content mutation (no-op text assignment leaves the\n    # bookmark untouched) so the run split collapses into a single run and\n    # the _GoBack bookmark that sat between the two runs is dropped.\n    $rng1.Text = \"{{ Name_2 }}~~TMP~~\"\n    $rng1.Text = \"{{ Name_2 }}\"\n}\n\n# --- Change 2: \"Firma del Dise\u00f1ador\" -> \"Firma de\" + _GoBack bookmark + \" Dise\u00f1o\" (bold) ---\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"Firma del Dise\u00f1ador\"\n$rng2.Find.Execute() | Out-Null\nif ($rng2.Find.Found) {\n    $startPos = $rng2.Start\n    $rng2.Text = \"Firma de Dise\u00f1o\"\n    $splitPos = $startPos + 8\n    $bmRange = $d.Range($splitPos, $splitPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
